# "Generate Report for Handback" — reflects that the handback files are now
# in sync with en-US: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", handback timestamps are refreshed, and
# the stale "handback file is not the latest" error on the a36f76a9 row is
# cleared now that it is back in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns widened to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn detail sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
# Latest Handback DateTime for the 3907f169 row refreshed.
$wsZh.Range("K2").Value = "2016-11-15 18:04:58"
# Stale "handback file is not the latest" error cleared for the a36f76a9 row.
$wsZh.Range("P3").Value = ""

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de detail sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
# Latest Handback DateTime for the 3907f169 row refreshed.
$wsDe.Range("K2").Value = "2016-11-15 18:05:18"
# Stale "handback file is not the latest" error cleared for the a36f76a9 row.
$wsDe.Range("P3").Value = ""

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(16).ColumnWidth = 12.833333333333334
